# Update the EDCR results sheet:
#  - Precision/Recall/F1 (Base Model) columns (H,I,J) are re-synced to the
#    Base Precision/Recall/F1 columns (E,F,G).
#  - The downstream Improvement / Improvement(%) columns (O,P,Q,R,S,T) are
#    recomputed from the refreshed base-model values.
#  - The "Label" column text is refreshed: the old combined
#    "Best Precision, Best Recall, Best F1" label is split into
#    "Best Recall, Best F1" (rows that already used it) and a new
#    "Best Precision" label (rows that used to read "Worst F1"); the
#    "Worst F1" label now also gets applied to the remaining (previously
#    unlabeled) rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $baseP = $ws.Cells.Item($r, 5).Value()   # E: Base Precision
    $baseR = $ws.Cells.Item($r, 6).Value()   # F: Base Recall
    $baseF = $ws.Cells.Item($r, 7).Value()   # G: Base F1

    # H/I/J: Precision/Recall/F1 (Base Model) now mirror Base Precision/Recall/F1
    $ws.Cells.Item($r, 8).Value = $baseP
    $ws.Cells.Item($r, 9).Value = $baseR
    $ws.Cells.Item($r, 10).Value = $baseF

    $edcrP = $ws.Cells.Item($r, 12).Value()  # L: Precision (EDCR)
    $edcrR = $ws.Cells.Item($r, 13).Value()  # M: Recall (EDCR)
    $edcrF = $ws.Cells.Item($r, 14).Value()  # N: F1 (EDCR)

    $impP = $edcrP - $baseP
    $impR = $edcrR - $baseR
    $impF = $edcrF - $baseF

    # O/P/Q: Precision/Recall/F1 Improvement
    $ws.Cells.Item($r, 15).Value = $impP
    $ws.Cells.Item($r, 16).Value = $impR
    $ws.Cells.Item($r, 17).Value = $impF

    if ($baseP -ne 0) { $impPPct = $impP / $baseP } else { $impPPct = 0 }
    if ($baseR -ne 0) { $impRPct = $impR / $baseR } else { $impRPct = 0 }
    if ($baseF -ne 0) { $impFPct = $impF / $baseF } else { $impFPct = 0 }

    # R/S/T: Precision/Recall/F1 Improvement (%)
    $ws.Cells.Item($r, 18).Value = $impPPct
    $ws.Cells.Item($r, 19).Value = $impRPct
    $ws.Cells.Item($r, 20).Value = $impFPct

    # U: Label text refresh
    if ($r -le 21) {
        $ws.Cells.Item($r, 21).Value = "Best Recall, Best F1"
    } elseif ($r -le 31) {
        $ws.Cells.Item($r, 21).Value = "Best Precision"
    } else {
        $ws.Cells.Item($r, 21).Value = "Worst F1"
    }
}
